# Updated symbol list (crypto price/volume refresh) — GitHub Actions run.
# Price cells in column D are stored as text (they look numeric, e.g. "247.87"),
# so each is written with a leading apostrophe to force text entry and then
# restyled back to "Normal" so no stray number-format/quote-prefix style sticks
# around (Excel auto-adds a style id for the quote-prefix flag otherwise).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'247.87"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'21.83"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Value = "'5.497"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = "'0.05653"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'3.380"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'6.436"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'0.8017"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'1.039"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Value = "'0.1426"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Value = "'0.07259"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "'0.03138"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "'0.02951"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = "'0.09281"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'0.001643"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "'3.218"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'0.04729"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "'0.0005853"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "17OneONEWorstin24h"
$ws.Range("D19").Value = "'0.006422"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Value = "'0.005025"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Value = "'0.001051"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = "'0.0001502"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Value = "'0.0003202"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Value = "'4.022"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Value = "'2.109"
$ws.Range("D25").Style = "Normal"
$ws.Range("D40").Value = "'0.04087"
$ws.Range("D40").Style = "Normal"
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").Value = "'0.1041"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "40BKEXTokenBKK"
$ws.Range("D42").Value = "'0.002973"
$ws.Range("D42").Style = "Normal"
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D43").Value = "'0.003269"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "42KickTokenKICK"
$ws.Range("D44").Value = "'0.009141"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005821"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Value = "'0.00000000750"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Value = "'0.7856"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Value = "'0.01661"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "47BOLOBOLO"
$ws.Range("D49").Value = "'0.00002101"
$ws.Range("D49").Style = "Normal"
